$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.199.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.477.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "488.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +12.22%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +4.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.487.62"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.48%  "
$ws.Range("E10").Value = "  +9.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0970"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.34%  "
$ws.Range("E12").Value = "  +5.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.124"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.910.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.159.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.55%  "
$ws.Range("E17").Value = "  +3.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.480.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.80%  "
$ws.Range("E19").Value = "  +8.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "318.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.92%  "
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.17%  "
$ws.Range("E25").Value = "  +7.18%  "
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("E27").Value = "  +3.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.585.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.90%  "
$ws.Range("E29").Value = "  +7.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0790"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "148.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.14%  "
$ws.Range("E34").Value = "  +5.49%  "
$ws.Range("E35").Value = "  +4.33%  "
$ws.Range("E36").Value = "  +8.54%  "
$ws.Range("E37").Value = "  +6.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.859"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.83%  "
$ws.Range("E40").Value = "  +8.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.995"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0554"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.605"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.99%  "
$ws.Range("E44").Value = "  +8.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0924"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "258.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +13.67%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0228"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.882.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.30%  "

Write-Output "done"